# Dividend Calculation.xlsx update
# - Taxable Account dividend for month 6 (row 8) on the "Yearly" sheet increased
#   from 118.95 to 142.2 (+23.25), which ripples through the dependent totals
#   (row 8's Total Dividends, the Taxable/Total column sums in row 15, and the
#   cross-sheet references + grand totals on the "All Time" sheet).
# - The active selection moves from Yearly!K22 to Yearly!H21, and the
#   workbook's active tab switches from "Yearly" to "All Time", landing on
#   All Time!N29 (scrolled so row 19 is at the top).

$wb = $excel.ActiveWorkbook

$yearly = $wb.Worksheets.Item("Yearly")
$allTime = $wb.Worksheets.Item("All Time")

# Update the Taxable Account figure for June 2017; dependent SUM/shared
# formulas (O8, L15, O15 on Yearly and F8, I8, F46, I46 on All Time)
# recalculate automatically.
$yearly.Range("L8").Value = 142.2

# Yearly sheet is active to start with - move its selection before switching.
$yearly.Range("H21").Select()

# Switch the active tab to "All Time" and update its view/selection.
$allTime.Activate()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$allTime.Range("N29").Select()
